$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C31").Value = 51
$ws.Range("E31").Value = 0.02205882352941177

$ws.Range("C34").Value = 84
$ws.Range("E34").Value = 0.03723404255319149

$ws.Range("C36").Value = 140
$ws.Range("E36").Value = 0.07253886010362694

$ws.Range("C37").Value = 895
$ws.Range("D37").Value = 895
